# The commit swaps the contents of ppt/theme/theme1.xml (the "Office Theme"
# color/font/format scheme) and ppt/theme/theme2.xml (the "Integral" color
# scheme that the deck's single slide master actually uses).
#
# theme1.xml and theme2.xml already share an identical fontScheme and
# fmtScheme - the only real difference between them is the <a:clrScheme>
# (and the cosmetic name="" attributes). So, for the part of the deck that
# is reachable through the PowerPoint object model (the slide master's
# theme, i.e. ppt/theme/theme2.xml), reproduce the swap by rewriting the
# 12 theme colors from the "Integral" palette to the "Office Theme"
# palette, via Slide.ThemeColorScheme (maps 1:1 onto
# dk1,lt1,dk2,lt2,accent1-6,hlink,folHlink).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

function ToRgbInt([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Target ("Office Theme") srgbClr values, in ThemeColorSchemeIndex order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3, 8 accent4,
# 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$officeThemeColors = @{
    1  = "000000"
    2  = "FFFFFF"
    3  = "44546A"
    4  = "E7E6E6"
    5  = "5B9BD5"
    6  = "ED7D31"
    7  = "A5A5A5"
    8  = "FFC000"
    9  = "4472C4"
    10 = "70AD47"
    11 = "0563C1"
    12 = "954F72"
}

for ($i = 1; $i -le 12; $i++) {
    $c = $tcs.Colors($i)
    $c.RGB = ToRgbInt($officeThemeColors[$i])
}
